$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header labels: "_old" -> "_FV2410", "_new" -> "_FV2504"
#    (columns A1:J1 are the FV2410 ["_old"] headers, K1 is "diff",
#     L1:U1 are the FV2504 ["_new"] headers)
# ---------------------------------------------------------------------------
$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
$fv2504Headers = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2410Headers[$i]
}
$ws.Cells.Item(1, 11).Value = "diff"
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2504Headers[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn the data range into an Excel Table ("Table1") spanning A1:U80
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U80")
$listObj = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$listObj.Name = "Table1"

# ---------------------------------------------------------------------------
# 3) Freeze the header row (row 1) on the sheet
# ---------------------------------------------------------------------------
$ws.Activate()
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
